$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.114.74"
$ws.Range("E2").Value = "'  -2.82%  "
$ws.Range("D3").Value = "'1.869.88"
$ws.Range("E3").Value = "'  -1.91%  "
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("D5").Value = "'307.47"
$ws.Range("E5").Value = "'  -1.89%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  +0.14%  "
$ws.Range("D7").Value = "'0.5064"
$ws.Range("E7").Value = "'  +1.56%  "
$ws.Range("D8").Value = "'0.3750"
$ws.Range("D9").Value = "'0.07149"
$ws.Range("E9").Value = "'  -1.92%  "
$ws.Range("D10").Value = "'0.8890"
$ws.Range("E10").Value = "'  -2.50%  "
$ws.Range("D11").Value = "'20.67"
$ws.Range("E11").Value = "'  -2.02%  "
$ws.Range("B12").Value = "'TRON"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07561"
$ws.Range("E12").Value = "'  -1.63%  "
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.856.44"
$ws.Range("E13").Value = "'  -2.74%  "
$ws.Range("D14").Value = "'5.322"
$ws.Range("E14").Value = "'  -3.25%  "
$ws.Range("D15").Value = "'89.32"
$ws.Range("E15").Value = "'  -3.23%  "
$ws.Range("E16").Value = "'  +0.19%  "
$ws.Range("D17").Value = "'0.000008491"
$ws.Range("E17").Value = "'  -2.74%  "
$ws.Range("E18").Value = "'  -3.46%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "'  +0.12%  "
$ws.Range("D20").Value = "'27.174.65"
$ws.Range("E20").Value = "'  -2.71%  "
$ws.Range("D21").Value = "'5.075"
$ws.Range("E21").Value = "'  -1.99%  "
$ws.Range("D22").Value = "'2.103.59"
$ws.Range("E22").Value = "'  -2.63%  "
$ws.Range("D23").Value = "'10.59"
$ws.Range("E23").Value = "'  -2.42%  "
$ws.Range("D24").Value = "'6.482"
$ws.Range("E24").Value = "'  -1.48%  "
$ws.Range("D25").Value = "'150.94"
$ws.Range("E25").Value = "'  -1.35%  "
$ws.Range("D26").Value = "'1.840"
$ws.Range("E26").Value = "'  -1.28%  "
$ws.Range("D27").Value = "'18.01"
$ws.Range("E27").Value = "'  -2.08%  "
$ws.Range("D28").Value = "'2.097"
$ws.Range("E28").Value = "'  -5.67%  "
$ws.Range("E29").Value = "'  -2.15%  "
$ws.Range("E30").Value = "'  -2.82%  "
$ws.Range("D31").Value = "'4.683"
$ws.Range("E31").Value = "'  -3.76%  "
$ws.Range("D32").Value = "'0.09015"
$ws.Range("E32").Value = "'  +0.04%  "
$ws.Range("D33").Value = "'0.05129"
$ws.Range("E33").Value = "'  -2.85%  "
$ws.Range("D34").Value = "'3.096"
$ws.Range("E34").Value = "'  -3.27%  "
$ws.Range("D35").Value = "'0.7421"
$ws.Range("E35").Value = "'  -4.02%  "
$ws.Range("D36").Value = "'1.160"
$ws.Range("E36").Value = "'  -5.95%  "
$ws.Range("D37").Value = "'0.02037"
$ws.Range("E37").Value = "'  -2.48%  "
$ws.Range("D38").Value = "'2.535"
$ws.Range("E38").Value = "'  -0.97%  "
$ws.Range("D39").Value = "'3.044"
$ws.Range("E39").Value = "'  -0.59%  "
$ws.Range("D40").Value = "'1.076"
$ws.Range("E40").Value = "'  -1.59%  "
$ws.Range("D41").Value = "'0.5387"
$ws.Range("E41").Value = "'  -2.97%  "
$ws.Range("D42").Value = "'6.585"
$ws.Range("E42").Value = "'  -4.42%  "
$ws.Range("D43").Value = "'115.31"
$ws.Range("E43").Value = "'  +2.13%  "
$ws.Range("D44").Value = "'8.427"
$ws.Range("E44").Value = "'  -0.93%  "
$ws.Range("D45").Value = "'0.1478"
$ws.Range("E45").Value = "'  -2.87%  "
$ws.Range("D46").Value = "'0.4643"
$ws.Range("E46").Value = "'  -4.09%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "'  +0.17%  "
$ws.Range("D48").Value = "'9.983"
$ws.Range("E48").Value = "'  -5.99%  "
$ws.Range("D49").Value = "'1.567"
$ws.Range("E49").Value = "'  -4.17%  "
$ws.Range("D50").Value = "'64.59"
$ws.Range("E50").Value = "'  -4.32%  "
$ws.Range("D51").Value = "'36.57"
$ws.Range("E51").Value = "'  -1.64%  "